$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new accelerometer/gyroscope reading was spliced into the top of this
# "falling" sample window: every existing reading (columns C:H) shifts down
# one row, a freshly captured reading becomes row 2, and the oldest reading
# (previously the last row, timestamp 2000) is dropped so the window keeps a
# constant sample count. The timestamp/label columns (A:B) are an unrelated
# fixed per-row sequence and are left untouched.

# ax, ay, az, gx, gy, gz for the new row 2 followed by every row that shifts
# down to take the place of the row above it (rows 3..21).
$readings = @(
    @(0.3787193298339841, 0.750096321105957, -2.375997304916381, -0.0930042341351509, 0.1044579595327377, -0.295353353023529),
    @(0.1436797380447386, 0.7403357923030853, -1.325827866792679, 0.0355829000473022, 0.3640756905078888, -0.06383541971445079),
    @(-1.716686248779299, 1.053612291812897, 0.4070562124252333, 0.5752823352813721, 0.5053382515907288, -0.007177666760981),
    @(-3.969705402851109, 0.09983259439468162, 0.9274015724658967, 0.5288565754890442, 0.7619016766548157, 0.2580905556678772),
    @(-5.239097833633425, -0.5959589481353762, -0.1143757104873694, -0.7357872128486633, -0.0797179117798805, 0.8136724829673767),
    @(2.036354780197156, -0.2977316975593562, -1.156379550695419, -1.146135926246643, -0.6188064813613892, -0.1640173196792602),
    @(5.741946458816498, -4.677600264549274, -8.871290028095274, 0.8213083148002625, -3.046079635620117, 1.09803032875061),
    @(-18.10554087162019, -17.49406802654266, -29.53683829307556, 1.518305540084839, -0.5958990454673767, 0.1006400510668754),
    @(1.688319206237781, -2.673514366149895, -12.34325218200682, -0.9292787313461304, 0.0415388382971286, -2.335643291473389),
    @(-3.507262408733383, 0.2684899270534595, -2.699394106864903, 0.5500841736793518, 1.703092336654663, -0.4928155243396759),
    @(-0.7152169942855791, 0.4368197321891781, 0.5279676914215112, 0.2722931802272796, -0.3020728528499603, 0.2756529450416565),
    @(0.77803122997284, 0.2878375947475432, -0.8067402243614243, -0.6982190012931824, 0.4492913782596588, -0.5900958180427551),
    @(-0.02435183525085627, -0.3931519985198989, -2.55875074863434, -0.4308127164840698, 1.22447943687439, -0.319024384021759),
    @(-0.2728092074394176, -0.8767854124307625, -2.626779749989506, -0.0684169083833694, 0.4977024495601654, -0.0363464802503585),
    @(1.066039085388185, -0.5432969331741325, -1.81298840045929, -0.1298088580369949, 0.0172569435089826, 0.1985312104225158),
    @(1.15899240970611, -0.2723855525255186, -1.873978555202483, 0.0598647929728031, 0.5253441333770752, -0.0415388382971286),
    @(-0.3004360198974618, 0.163273096084592, -1.427715420722963, -0.1009454801678657, 0.2724458873271942, -0.4051563739776611),
    @(-0.5053797960281373, -0.4524855315685276, -1.996987149119378, -0.0010690141934901, -0.1505782902240753, 0.0145080499351024),
    @(0.1984210014343271, -0.3136261999607083, -1.604727536439895, -0.1108720451593399, 0.0652098655700683, -0.131183311343193),
    @(0.1348390579223627, -0.3085210472345353, -1.704802349209786, 0.0563523173332214, 0.0316122770309448, 0.1798998117446899)
)

for ($i = 0; $i -lt $readings.Length; $i++) {
    $r = $i + 2
    $reading = $readings[$i]
    $ws.Range("C$r").Value = $reading[0]
    $ws.Range("D$r").Value = $reading[1]
    $ws.Range("E$r").Value = $reading[2]
    $ws.Range("F$r").Value = $reading[3]
    $ws.Range("G$r").Value = $reading[4]
    $ws.Range("H$r").Value = $reading[5]
}

# Drop the now-duplicated last row (old row 22 / timestamp 2000) and shrink
# the used range back down to 21 rows.
$ws.Rows.Item(22).Delete()
